# Fruta / hortaliza, semanal
# Insert a new weekly price-report row (row 37) into the Melón subset sheet,
# pushing all subsequent rows down by one (old row 37 -> 38, ..., old row 82 -> 83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 37:82 down to 38:83, creating a blank row 37.
$ws.Rows(37).Insert()

# Populate the new row 37 with the latest weekly measurement.
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 45012
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112027
$ws.Range("G37").Value = "Melón"
$ws.Range("H37").Value = "Tuna"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 140
$ws.Range("K37").Value = 17000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17429
$ws.Range("N37").Value = "$/caja 18 unidades"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 968
$ws.Range("Q37").Value = 18
$ws.Range("R37").Value = "Hortaliza"
